# Update the regression-coefficient table (adds "Crisis and Credit
# Allocation" period results) by replacing the 9 data-cell values in
# B2:D4 with their new figures, while keeping the cells as plain text
# (matching the original shared-string / "t=s" storage) rather than
# letting Excel auto-coerce the purely-numeric-looking ones into
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "bare" number (no trailing significance
# stars) need to be pre-formatted as Text, otherwise Excel's COM layer
# will silently convert the assigned string into a numeric cell.
$numericLookingCells = @("B2", "B3", "B4", "C4", "D2")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("B2").Value = "0.17"
$ws.Range("B3").Value = "-0.01"
$ws.Range("B4").Value = "-0.09"

$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
$ws.Range("C4").Value = "0.98"

$ws.Range("D2").Value = "-0.89"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"
